$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 232.75
$ws.Range("J53").Value = 115.5
$ws.Range("L53").Value = 115.5
$ws.Range("N53").Value = -1389.5
$ws.Range("H88").Value = 3000
$ws.Range("J88").Value = 3000
$ws.Range("L88").Value = 3000
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 3000
$ws.Range("J91").Value = 3000
$ws.Range("L91").Value = 3000
$ws.Range("N91").Value = -5808
$ws.Range("H99").Value = 287
$ws.Range("I99").Value = 287
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 861
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 637
$ws.Range("N99").ClearContents()
$ws.Range("H106").Value = 14797.8
$ws.Range("I106").Value = 13497.25
$ws.Range("J106").Value = 20000
$ws.Range("K106").Value = 13497.25
$ws.Range("L106").Value = 20000
$ws.Range("M106").Value = -12866.25
$ws.Range("N106").Value = -21262
$ws.Range("H132").Value = 1325
$ws.Range("I132").Value = 1341.2142
$ws.Range("K132").Value = 4023.6426
$ws.Range("M132").Value = -1493.6426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 7966.6665
$ws.Range("I16").Value = 14000
$ws.Range("J16").Value = 4950
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 4950
$ws.Range("M16").Value = -13713
$ws.Range("N16").Value = -5524
$ws.Range("H32").Value = 4615.1113
$ws.Range("I32").Value = 3183.468
$ws.Range("K32").Value = 3183.468
$ws.Range("M32").Value = -2896.468
$ws.Range("H61").Value = 2075.2942
$ws.Range("I61").Value = 2037.9333
$ws.Range("K61").Value = 2037.9333
$ws.Range("M61").Value = -1825.9333
$ws.Range("H74").Value = 531.8788
$ws.Range("I74").Value = 531.8788
$ws.Range("K74").Value = 531.8788
$ws.Range("M74").Value = 342.1212
$ws.Range("H77").Value = 531.8788
$ws.Range("I77").Value = 531.8788
$ws.Range("K77").Value = 2659.394
$ws.Range("M77").Value = 1708.606
$ws.Range("H136").Value = 2075.2942
$ws.Range("I136").Value = 2037.9333
$ws.Range("K136").Value = 6113.7999
$ws.Range("M136").Value = -3563.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 145.16667
$ws.Range("I22").Value = 145.16667
$ws.Range("K22").Value = 145.16667
$ws.Range("M22").Value = 27.83332999999999
$ws.Range("H64").Value = 1931.3334
$ws.Range("I64").Value = 1700
$ws.Range("J64").Value = 1977.6
$ws.Range("K64").Value = 1700
$ws.Range("L64").Value = 1977.6
$ws.Range("N64").Value = -2427.6
$ws.Range("M64").Value = -1475
$ws.Range("H67").Value = 1931.3334
$ws.Range("I67").Value = 1700
$ws.Range("J67").Value = 1977.6
$ws.Range("K67").Value = 1700
$ws.Range("L67").Value = 1977.6
$ws.Range("N67").Value = -3537.6
$ws.Range("M67").Value = -920
$ws.Range("H134").Value = 3046.75
$ws.Range("I134").Value = 2070.125
$ws.Range("K134").Value = 6210.375
$ws.Range("M134").Value = -3675.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H4").Value = 4243.75
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 4243.75
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 4243.75
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -4467.75
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H29").Value = 11620
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 11620
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 11620
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -12206
$ws.Range("H31").Value = 3643.122
$ws.Range("I31").Value = 2216.7144
$ws.Range("J31").Value = 6715.385
$ws.Range("K31").Value = 2216.7144
$ws.Range("L31").Value = 6715.385
$ws.Range("M31").Value = -1921.7144
$ws.Range("N31").Value = -7305.385
$ws.Range("H34").Value = 3643.122
$ws.Range("I34").Value = 2216.7144
$ws.Range("J34").Value = 6715.385
$ws.Range("K34").Value = 2216.7144
$ws.Range("L34").Value = 6715.385
$ws.Range("M34").Value = -2014.7144
$ws.Range("N34").Value = -7119.385
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H134").Value = 2332.2896
$ws.Range("I134").Value = 2209.7778
$ws.Range("J134").Value = 2633
$ws.Range("K134").Value = 6629.3334
$ws.Range("L134").Value = 7899
$ws.Range("M134").Value = -4094.3334
$ws.Range("N134").Value = -12969

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2232.2666
$ws.Range("I126").Value = 2320.875
$ws.Range("J126").Value = 2131
$ws.Range("K126").Value = 6962.625
$ws.Range("L126").Value = 6393
$ws.Range("M126").Value = -4492.625
$ws.Range("N126").Value = -11333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 9498
$ws.Range("J3").Value = 9498
$ws.Range("L3").Value = 9498
$ws.Range("N3").Value = -9722
$ws.Range("H4").Value = 3875
$ws.Range("I4").Value = 5000
$ws.Range("J4").Value = 2750
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 2750
$ws.Range("M4").Value = -4887
$ws.Range("N4").Value = -2976
$ws.Range("H5").Value = 15011
$ws.Range("J5").Value = 15011
$ws.Range("L5").Value = 15011
$ws.Range("N5").Value = -15237
$ws.Range("H15").Value = 9498
$ws.Range("J15").Value = 9498
$ws.Range("L15").Value = 9498
$ws.Range("N15").Value = -9838
$ws.Range("H18").Value = 35000
$ws.Range("I18").Value = 35000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 35000
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("M18").Value = -34828
$ws.Range("H23").Value = 10502500
$ws.Range("I23").Value = 21000000
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 21000000
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = -20999770
$ws.Range("N23").Value = -5460
$ws.Range("H28").Value = 3875
$ws.Range("I28").Value = 5000
$ws.Range("J28").Value = 2750
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 2750
$ws.Range("M28").Value = -4768
$ws.Range("N28").Value = -3214
$ws.Range("H34").Value = 19000
$ws.Range("I34").Value = 30000
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 30000
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -29828
$ws.Range("N34").Value = -8344
$ws.Range("H37").Value = 3875
$ws.Range("I37").Value = 5000
$ws.Range("J37").Value = 2750
$ws.Range("K37").Value = 5000
$ws.Range("L37").Value = 2750
$ws.Range("M37").Value = -4893
$ws.Range("N37").Value = -2964
$ws.Range("H43").Value = 8435.5
$ws.Range("J43").Value = 8469.076999999999
$ws.Range("L43").Value = 8469.076999999999
$ws.Range("N43").Value = -8855.076999999999
$ws.Range("H136").Value = 5411.9473
$ws.Range("I136").Value = 5349.353
$ws.Range("J136").Value = 5944
$ws.Range("K136").Value = 16048.059
$ws.Range("L136").Value = 17832
$ws.Range("M136").Value = -13498.059
$ws.Range("N136").Value = -22932

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6500
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 6500
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 6500
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -6780
$ws.Range("H11").Value = 17502
$ws.Range("I11").Value = 14999
$ws.Range("J11").Value = 20005
$ws.Range("K11").Value = 14999
$ws.Range("L11").Value = 20005
$ws.Range("M11").Value = -14857
$ws.Range("N11").Value = -20289
$ws.Range("H81").Value = 3294.9
$ws.Range("J81").Value = 1000.6667
$ws.Range("L81").Value = 2001.3334
$ws.Range("N81").Value = -4123.3334
$ws.Range("H84").Value = 3294.9
$ws.Range("J84").Value = 1000.6667
$ws.Range("L84").Value = 10006.667
$ws.Range("N84").Value = -20614.667
$ws.Range("H118").Value = 39666.668
$ws.Range("J118").Value = 39666.668
$ws.Range("L118").Value = 39666.668
$ws.Range("N118").Value = -42980.668
$ws.Range("H132").Value = 3313.5
$ws.Range("I132").Value = 2626.3333
$ws.Range("K132").Value = 7878.999899999999
$ws.Range("M132").Value = -5348.999899999999
$ws.Range("H136").Value = 1125.0938
$ws.Range("I136").Value = 714.4286
$ws.Range("K136").Value = 2143.2858
$ws.Range("M136").Value = 406.7142000000003
